$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2,3).Value = 5
$ws.Cells.Item(2,6).Value = 269
$ws.Cells.Item(2,8).Value = "bedrooms"
$ws.Cells.Item(2,9).Value = "distractor"
$ws.Cells.Item(2,11).Value = "f"
$ws.Cells.Item(2,12).Value = "stimuli/img_n9xll.png"
$ws.Cells.Item(2,13).Value = 77.14285714285714
$ws.Cells.Item(2,14).Value = 59.21428571428572
$ws.Cells.Item(2,15).Value = 68.17857142857143
$ws.Cells.Item(2,16).Value = 42
$ws.Cells.Item(2,17).Value = 7
$ws.Cells.Item(2,18).Value = 7
$ws.Cells.Item(2,19).Value = 7

# Row 3
$ws.Cells.Item(3,3).Value = 5
$ws.Cells.Item(3,6).Value = 270
$ws.Cells.Item(3,12).Value = "stimuli/img_60242.png"
$ws.Cells.Item(3,13).Value = 78.33333333333333
$ws.Cells.Item(3,14).Value = 57.57575757575758
$ws.Cells.Item(3,15).Value = 67.95454545454545
$ws.Cells.Item(3,16).Value = 33
$ws.Cells.Item(3,17).Value = 7
$ws.Cells.Item(3,18).Value = 7
$ws.Cells.Item(3,19).Value = 7

# Row 4
$ws.Cells.Item(4,3).Value = 5
$ws.Cells.Item(4,6).Value = 271
$ws.Cells.Item(4,8).Value = "kitchens"
$ws.Cells.Item(4,9).Value = "target"
$ws.Cells.Item(4,11).Value = "j"
$ws.Cells.Item(4,12).Value = "stimuli/img_inqod.png"
$ws.Cells.Item(4,13).Value = 70.84848484848484
$ws.Cells.Item(4,14).Value = 50.63636363636363
$ws.Cells.Item(4,15).Value = 60.74242424242424
$ws.Cells.Item(4,16).Value = 33
$ws.Cells.Item(4,17).Value = 5
$ws.Cells.Item(4,18).Value = 5
$ws.Cells.Item(4,19).Value = 5

# Row 5
$ws.Cells.Item(5,3).Value = 5
$ws.Cells.Item(5,6).Value = 272
$ws.Cells.Item(5,8).Value = "kitchens"
$ws.Cells.Item(5,9).Value = "target"
$ws.Cells.Item(5,11).Value = "j"
$ws.Cells.Item(5,12).Value = "stimuli/img_9mky8.png"
$ws.Cells.Item(5,13).Value = 84.32352941176471
$ws.Cells.Item(5,14).Value = 65.17647058823529
$ws.Cells.Item(5,15).Value = 74.75
$ws.Cells.Item(5,16).Value = 34
$ws.Cells.Item(5,17).Value = 9
$ws.Cells.Item(5,18).Value = 9
$ws.Cells.Item(5,19).Value = 9

# Row 6
$ws.Cells.Item(6,3).Value = 5
$ws.Cells.Item(6,6).Value = 273
$ws.Cells.Item(6,12).Value = "stimuli/img_cv6mf.png"
$ws.Cells.Item(6,13).Value = 66.8
$ws.Cells.Item(6,14).Value = 42.08
$ws.Cells.Item(6,15).Value = 54.44
$ws.Cells.Item(6,16).Value = 25
$ws.Cells.Item(6,17).Value = 4
$ws.Cells.Item(6,18).Value = 4
$ws.Cells.Item(6,19).Value = 4

# Row 7
$ws.Cells.Item(7,3).Value = 5
$ws.Cells.Item(7,6).Value = 274
$ws.Cells.Item(7,12).Value = "stimuli/img_iyxnj.png"
$ws.Cells.Item(7,13).Value = 75.30555555555556
$ws.Cells.Item(7,14).Value = 54.33333333333334
$ws.Cells.Item(7,15).Value = 64.81944444444444
$ws.Cells.Item(7,16).Value = 36

# Row 8
$ws.Cells.Item(8,3).Value = 5
$ws.Cells.Item(8,6).Value = 275
$ws.Cells.Item(8,12).Value = "stimuli/img_wyl6z.png"
$ws.Cells.Item(8,13).Value = 59.8235294117647
$ws.Cells.Item(8,14).Value = 36.23529411764706
$ws.Cells.Item(8,15).Value = 48.02941176470588
$ws.Cells.Item(8,16).Value = 34
$ws.Cells.Item(8,17).Value = 3
$ws.Cells.Item(8,18).Value = 3
$ws.Cells.Item(8,19).Value = 3

# Row 9
$ws.Cells.Item(9,3).Value = 5
$ws.Cells.Item(9,6).Value = 276
$ws.Cells.Item(9,8).Value = "living_rooms"
$ws.Cells.Item(9,9).Value = "distractor"
$ws.Cells.Item(9,11).Value = "f"
$ws.Cells.Item(9,12).Value = "stimuli/img_hmmra.png"
$ws.Cells.Item(9,13).Value = 54.65853658536585
$ws.Cells.Item(9,14).Value = 34.24390243902439
$ws.Cells.Item(9,15).Value = 44.45121951219512
$ws.Cells.Item(9,16).Value = 41
$ws.Cells.Item(9,17).Value = 3
$ws.Cells.Item(9,18).Value = 3
$ws.Cells.Item(9,19).Value = 3

# Row 10
$ws.Cells.Item(10,3).Value = 5
$ws.Cells.Item(10,6).Value = 277
$ws.Cells.Item(10,8).Value = "kitchens"
$ws.Cells.Item(10,9).Value = "target"
$ws.Cells.Item(10,11).Value = "j"
$ws.Cells.Item(10,12).Value = "stimuli/img_eatdk.png"
$ws.Cells.Item(10,13).Value = 81.40625
$ws.Cells.Item(10,14).Value = 61.375
$ws.Cells.Item(10,15).Value = 71.390625
$ws.Cells.Item(10,16).Value = 32

# Row 11
$ws.Cells.Item(11,3).Value = 5
$ws.Cells.Item(11,6).Value = 278
$ws.Cells.Item(11,12).Value = "stimuli/img_uwv6y.png"
$ws.Cells.Item(11,13).Value = 78.88888888888889
$ws.Cells.Item(11,14).Value = 59.30555555555556
$ws.Cells.Item(11,15).Value = 69.09722222222223
$ws.Cells.Item(11,16).Value = 36
$ws.Cells.Item(11,17).Value = 8
$ws.Cells.Item(11,18).Value = 8
$ws.Cells.Item(11,19).Value = 8

# Row 12
$ws.Cells.Item(12,3).Value = 5
$ws.Cells.Item(12,6).Value = 279
$ws.Cells.Item(12,12).Value = "stimuli/img_7wul8.png"
$ws.Cells.Item(12,13).Value = 43.03030303030303
$ws.Cells.Item(12,14).Value = 25.54545454545455
$ws.Cells.Item(12,15).Value = 34.28787878787879
$ws.Cells.Item(12,16).Value = 33

# Row 13
$ws.Cells.Item(13,3).Value = 5
$ws.Cells.Item(13,6).Value = 280
$ws.Cells.Item(13,12).Value = "stimuli/img_aplao.png"
$ws.Cells.Item(13,13).Value = 64.09090909090909
$ws.Cells.Item(13,14).Value = 40.75757575757576
$ws.Cells.Item(13,15).Value = 52.42424242424242
$ws.Cells.Item(13,17).Value = 3
$ws.Cells.Item(13,18).Value = 3
$ws.Cells.Item(13,19).Value = 3

# Row 14
$ws.Cells.Item(14,3).Value = 5
$ws.Cells.Item(14,6).Value = 281
$ws.Cells.Item(14,12).Value = "stimuli/img_79b5l.png"
$ws.Cells.Item(14,13).Value = 72.74285714285715
$ws.Cells.Item(14,14).Value = 53.31428571428572
$ws.Cells.Item(14,15).Value = 63.02857142857143
$ws.Cells.Item(14,16).Value = 35
$ws.Cells.Item(14,17).Value = 6
$ws.Cells.Item(14,18).Value = 6
$ws.Cells.Item(14,19).Value = 6

# Row 15
$ws.Cells.Item(15,3).Value = 5
$ws.Cells.Item(15,6).Value = 282
$ws.Cells.Item(15,12).Value = "stimuli/img_q1ynd.png"
$ws.Cells.Item(15,13).Value = 70.05714285714286
$ws.Cells.Item(15,14).Value = 47.31428571428572
$ws.Cells.Item(15,15).Value = 58.68571428571429
$ws.Cells.Item(15,17).Value = 5
$ws.Cells.Item(15,18).Value = 5
$ws.Cells.Item(15,19).Value = 5

# Row 16
$ws.Cells.Item(16,3).Value = 5
$ws.Cells.Item(16,6).Value = 283
$ws.Cells.Item(16,8).Value = "living_rooms"
$ws.Cells.Item(16,9).Value = "distractor"
$ws.Cells.Item(16,11).Value = "f"
$ws.Cells.Item(16,12).Value = "stimuli/img_3sw8t.png"
$ws.Cells.Item(16,13).Value = 67.48888888888889
$ws.Cells.Item(16,14).Value = 48.51111111111111
$ws.Cells.Item(16,15).Value = 58
$ws.Cells.Item(16,16).Value = 45
$ws.Cells.Item(16,17).Value = 5
$ws.Cells.Item(16,18).Value = 5
$ws.Cells.Item(16,19).Value = 5

# Row 17
$ws.Cells.Item(17,3).Value = 5
$ws.Cells.Item(17,6).Value = 284
$ws.Cells.Item(17,8).Value = "living_rooms"
$ws.Cells.Item(17,9).Value = "distractor"
$ws.Cells.Item(17,11).Value = "f"
$ws.Cells.Item(17,12).Value = "stimuli/img_kq9s9.png"
$ws.Cells.Item(17,13).Value = 62.30232558139535
$ws.Cells.Item(17,14).Value = 39.97674418604651
$ws.Cells.Item(17,15).Value = 51.13953488372093
$ws.Cells.Item(17,16).Value = 43
$ws.Cells.Item(17,17).Value = 4
$ws.Cells.Item(17,18).Value = 4
$ws.Cells.Item(17,19).Value = 4

# Row 18
$ws.Cells.Item(18,3).Value = 5
$ws.Cells.Item(18,6).Value = 285
$ws.Cells.Item(18,8).Value = "kitchens"
$ws.Cells.Item(18,9).Value = "target"
$ws.Cells.Item(18,11).Value = "j"
$ws.Cells.Item(18,12).Value = "stimuli/img_nyv2b.png"
$ws.Cells.Item(18,13).Value = 11.91176470588235
$ws.Cells.Item(18,14).Value = 6.852941176470588
$ws.Cells.Item(18,15).Value = 9.382352941176471
$ws.Cells.Item(18,16).Value = 34
$ws.Cells.Item(18,17).Value = 1
$ws.Cells.Item(18,18).Value = 1
$ws.Cells.Item(18,19).Value = 1

# Row 19
$ws.Cells.Item(19,3).Value = 5
$ws.Cells.Item(19,6).Value = 286
$ws.Cells.Item(19,8).Value = "kitchens"
$ws.Cells.Item(19,9).Value = "target"
$ws.Cells.Item(19,11).Value = "j"
$ws.Cells.Item(19,12).Value = "stimuli/img_d8xbu.png"
$ws.Cells.Item(19,13).Value = 91.36363636363636
$ws.Cells.Item(19,14).Value = 73.18181818181819
$ws.Cells.Item(19,15).Value = 82.27272727272728
$ws.Cells.Item(19,16).Value = 33

# Row 20
$ws.Cells.Item(20,3).Value = 5
$ws.Cells.Item(20,6).Value = 287
$ws.Cells.Item(20,8).Value = "living_rooms"
$ws.Cells.Item(20,9).Value = "distractor"
$ws.Cells.Item(20,11).Value = "f"
$ws.Cells.Item(20,12).Value = "stimuli/img_37hgm.png"
$ws.Cells.Item(20,13).Value = 70.95454545454545
$ws.Cells.Item(20,14).Value = 54.77272727272727
$ws.Cells.Item(20,15).Value = 62.86363636363636
$ws.Cells.Item(20,16).Value = 44
$ws.Cells.Item(20,17).Value = 6
$ws.Cells.Item(20,18).Value = 6
$ws.Cells.Item(20,19).Value = 6

# Row 21
$ws.Cells.Item(21,3).Value = 5
$ws.Cells.Item(21,6).Value = 288
$ws.Cells.Item(21,8).Value = "living_rooms"
$ws.Cells.Item(21,9).Value = "distractor"
$ws.Cells.Item(21,11).Value = "f"
$ws.Cells.Item(21,12).Value = "stimuli/img_iudc4.png"
$ws.Cells.Item(21,13).Value = 73.625
$ws.Cells.Item(21,14).Value = 52.275
$ws.Cells.Item(21,15).Value = 62.95
$ws.Cells.Item(21,16).Value = 40
$ws.Cells.Item(21,17).Value = 6
$ws.Cells.Item(21,18).Value = 6
$ws.Cells.Item(21,19).Value = 6

# Row 22
$ws.Cells.Item(22,3).Value = 5
$ws.Cells.Item(22,6).Value = 289
$ws.Cells.Item(22,8).Value = "kitchens"
$ws.Cells.Item(22,9).Value = "target"
$ws.Cells.Item(22,11).Value = "j"
$ws.Cells.Item(22,12).Value = "stimuli/img_ye5sl.png"
$ws.Cells.Item(22,13).Value = 53.2258064516129
$ws.Cells.Item(22,14).Value = 34.45161290322581
$ws.Cells.Item(22,15).Value = 43.83870967741936
$ws.Cells.Item(22,16).Value = 31
$ws.Cells.Item(22,17).Value = 2
$ws.Cells.Item(22,18).Value = 2
$ws.Cells.Item(22,19).Value = 2

# Row 23
$ws.Cells.Item(23,3).Value = 5
$ws.Cells.Item(23,6).Value = 290
$ws.Cells.Item(23,12).Value = "stimuli/img_ce9vx.png"
$ws.Cells.Item(23,13).Value = 75.90909090909091
$ws.Cells.Item(23,14).Value = 57.12121212121212
$ws.Cells.Item(23,15).Value = 66.51515151515152
$ws.Cells.Item(23,16).Value = 33
$ws.Cells.Item(23,17).Value = 7
$ws.Cells.Item(23,18).Value = 7
$ws.Cells.Item(23,19).Value = 7

# Row 24
$ws.Cells.Item(24,3).Value = 5
$ws.Cells.Item(24,6).Value = 291
$ws.Cells.Item(24,12).Value = "stimuli/img_a8wvq.png"
$ws.Cells.Item(24,13).Value = 86.25925925925925
$ws.Cells.Item(24,14).Value = 66.25925925925925
$ws.Cells.Item(24,15).Value = 76.25925925925925
$ws.Cells.Item(24,16).Value = 27
$ws.Cells.Item(24,17).Value = 10
$ws.Cells.Item(24,18).Value = 10
$ws.Cells.Item(24,19).Value = 10

# Row 25
$ws.Cells.Item(25,3).Value = 5
$ws.Cells.Item(25,6).Value = 292
$ws.Cells.Item(25,12).Value = "stimuli/img_yeh72.png"
$ws.Cells.Item(25,13).Value = 68.66666666666667
$ws.Cells.Item(25,14).Value = 45.21212121212121
$ws.Cells.Item(25,15).Value = 56.93939393939394
$ws.Cells.Item(25,16).Value = 33
$ws.Cells.Item(25,17).Value = 4
$ws.Cells.Item(25,18).Value = 4
$ws.Cells.Item(25,19).Value = 4

# Row 26
$ws.Cells.Item(26,3).Value = 5
$ws.Cells.Item(26,6).Value = 293
$ws.Cells.Item(26,12).Value = "stimuli/img_es7o2.png"
$ws.Cells.Item(26,13).Value = 52.48571428571429
$ws.Cells.Item(26,14).Value = 27.54285714285714
$ws.Cells.Item(26,15).Value = 40.01428571428572
$ws.Cells.Item(26,16).Value = 35

# Row 27
$ws.Cells.Item(27,3).Value = 5
$ws.Cells.Item(27,6).Value = 294
$ws.Cells.Item(27,12).Value = "stimuli/img_t90e2.png"
$ws.Cells.Item(27,13).Value = 83.0625
$ws.Cells.Item(27,14).Value = 61.96875
$ws.Cells.Item(27,15).Value = 72.515625
$ws.Cells.Item(27,16).Value = 32
